$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C110").Value = $null

$data = @(
    @(0, 8.308481844973947, $null),
    @(0, 8.308481844973947, $null),
    @(0, 8.308481844973947, 0.08425404025794732),
    @(0, 8.974591905842381, 0.4632327957854892),
    @(0, 8.974591905842381, 0.4646062144118369),
    @(0, 8.107829945440615, 0.1754485456778251),
    @(0, 8.039306644322307, 0.8160774729312702),
    @(0, 9.898999999999999, 0.4042252753665664)
)

$row = 111
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    if ($null -ne $r[2]) {
        $ws.Cells.Item($row, 3).Value = $r[2]
    }
    $row = $row + 1
}
